$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New titration/CRM accuracy reading taken 2021-04-18, appended as row 30
# (mirrors the existing row layout: Date, CRM value, Batch value, % off, Batch #, Notes)
$ws.Range("A30").Value = 20210418
$ws.Range("B30").Value = 2218.2559999999999
$ws.Range("C30").Value = 2224.4699999999998
$ws.Range("D30").Formula = "=100*(B30-C30)/C30"
$ws.Range("E30").Value = 180
$ws.Range("F30").Value = "CRM opened 20210418"

# Match the author's final selection/scroll position after adding the row
[void]$ws.Range("H30").Select()

$wb.Save()
